$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the "Invalid" (G) and "Absent" (H) counts for each attendance date row
# from 0 to 1 as recorded in the updated attendance sheet.

$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

$ws.Range("H5").Value = 1
$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1

$ws.Range("G11").Value = 1
$ws.Range("H11").Value = 1

$ws.Range("H12").Value = 1
$ws.Range("H13").Value = 1
$ws.Range("H14").Value = 1
$ws.Range("H15").Value = 1
$ws.Range("H16").Value = 1
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
